# This script normalizes vaccine/product label text across all worksheets:
#  - strips the trailing footnote markers like " [1]", " [2]", " [3]", " [4]", " [5]"
#    (replacing them with a single trailing space, matching the source text's
#    xml:space="preserve" trailing-space convention)
#  - collapses embedded line breaks in multi-line labels (e.g. "Recombivax\nHB",
#    "Fluzone\nQuadrivalent", "Influenza [5]\n(Age ...)") into single-line text
#    joined by a space
# The lookup table below maps each exact "before" cell text to its "after" text.

$wb = $excel.ActiveWorkbook

$map = @{}
$map["DTaP [1]"] = "DTaP "
$map["DTaP-IPV [2]"] = "DTaP-IPV "
$map["DTaP-Hep B-IPV [4]"] = "DTaP-Hep B-IPV "
$map["DTaP-IP-HI [4]"] = "DTaP-IP-HI "
$map["e-IPV [5]"] = "e-IPV "
$map["Hepatitis A Pediatric [5]"] = "Hepatitis A Pediatric "
$map["Hepatitis A-Hepatitis B 18 only [3]"] = "Hepatitis A-Hepatitis B 18 only "
$map["Hepatitis B [5]`nPediatric/Adolescent"] = "Hepatitis B  Pediatric/Adolescent"
$map["Recombivax`nHB"] = "Recombivax HB"
$map["Hib [5]"] = "Hib "
$map["HPV - Human Papillomavirus 9-valent [5]"] = "HPV - Human Papillomavirus 9-valent "
$map["MENB - Meningococcal Group B [5]"] = "MENB - Meningococcal Group B "
$map["Meningococcal Conjugate (Groups A, C, Y and W-135) [5]"] = "Meningococcal Conjugate (Groups A, C, Y and W-135) "
$map["Measles, Mumps and Rubella (MMR) [1]"] = "Measles, Mumps and Rubella (MMR) "
$map["MMR/Varicella [2]"] = "MMR/Varicella "
$map["Pneumococcal`n13-valent [5] (Pediatric)"] = "Pneumococcal 13-valent  (Pediatric)"
$map["Rotavirus, Live, Oral, Pentavalent [5]"] = "Rotavirus, Live, Oral, Pentavalent "
$map["Rotavirus, Live, Oral, Oral [5]"] = "Rotavirus, Live, Oral, Oral "
$map["Tetanus and Diphtheria Toxoids [3]"] = "Tetanus and Diphtheria Toxoids "
$map["Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis [1]"] = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis "
$map["Varicella [5]"] = "Varicella "
$map["Hepatitis A-Adult [5]"] = "Hepatitis A-Adult "
$map["Hepatitis A Adult [5]"] = "Hepatitis A Adult "
$map["Hepatitis A-Hepatitis B Adult [3]"] = "Hepatitis A-Hepatitis B Adult "
$map["Hepatitis B-Adult [5]"] = "Hepatitis B-Adult "
$map["HPV-Human Papillomavirus 9 Valent [5]"] = "HPV-Human Papillomavirus 9 Valent "
$map["Measles, Mumps,  Rubella-Adult [1]"] = "Measles, Mumps,  Rubella-Adult "
$map["Meningococcal Conjugate [5]"] = "Meningococcal Conjugate "
$map["Pneumococcal`n13-valent [5] (Adult)"] = "Pneumococcal 13-valent  (Adult)"
$map["Varicella-Adult [5]"] = "Varicella-Adult "
$map["Influenza [5]`n(Age 6 months and older)"] = "Influenza  (Age 6 months and older)"
$map["Fluzone`nQuadrivalent"] = "Fluzone Quadrivalent"
$map["Influenza [5]`n(Age 6-35 months)"] = "Influenza  (Age 6-35 months)"
$map["Fluzone`nQuadrivalent`nPediatric dose"] = "Fluzone Quadrivalent Pediatric dose"
$map["Influenza [5]`n(Age 36 months and older)"] = "Influenza  (Age 36 months and older)"
$map["Fluarix`nQuadrivalent"] = "Fluarix Quadrivalent"
$map["FluLaval`nQuadrivalent"] = "FluLaval Quadrivalent"
$map["Influenza [5]`n(Age 4 years and older)"] = "Influenza  (Age 4 years and older)"
$map["Influenza [5]`n(Age 9 years and older)"] = "Influenza  (Age 9 years and older)"
$map["10 pack-1 dose`nsyringe"] = "10 pack-1 dose syringe"
$map["Influenza [5]`n(Age 18 years and older)"] = "Influenza  (Age 18 years and older)"
$map["Afluria`nQuadrivalent"] = "Afluria Quadrivalent"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $startRow = $used.Row
    $startCol = $used.Column
    $numRows = $used.Rows.Count
    $numCols = $used.Columns.Count
    $endRow = $startRow + $numRows - 1
    $endCol = $startCol + $numCols - 1

    for ($r = $startRow; $r -le $endRow; $r++) {
        for ($c = $startCol; $c -le $endCol; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            $val = $cell.Value()
            if ($null -ne $val -and $val -is [string] -and $map.ContainsKey($val)) {
                $cell.Value = $map[$val]
            }
        }
    }
}
